$d = $word.ActiveDocument

# Insert the contact-info paragraph right after the "Dheeraj Chand" title
# paragraph, matching the centered-but-unstyled formatting used in the
# target document. Using Find/Replace with a ^p (paragraph mark) code lets
# the new paragraph/run pick up only the paragraph's alignment (jc=center)
# without inheriting the bold/large-size direct character formatting that
# a plain InsertParagraphAfter() would clone from the title run.
$d.Content.Find.Execute(
    "Dheeraj Chand",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
) | Out-Null
